$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column S, shifting old S:V to T:W.
# Excel copies the left-neighbour's formatting into the new column, so we
# also mirror column R's width onto the freshly inserted column S.
$ws.Columns("S").Insert() | Out-Null
$ws.Columns("S").ColumnWidth = $ws.Columns("R").ColumnWidth

# Row 2's new S2 cell needs the same (bordered) cell style the rest of the
# header row uses; copy formats from its row-2 neighbour so it reuses the
# existing style instead of minting a new one.
$ws.Range("T2").Copy() | Out-Null
$ws.Range("S2").PasteSpecial(-4122) | Out-Null

# Populate the new "saveProperties" command column (header + 2 argument rows).
$ws.Range("S1").Value = "saveProperties"
$ws.Range("S2").Value = "C:\Users\xihu_\Desktop\3.json"
$ws.Range("S3").Value = "name,password,age,password1"

# Update the active sheet view to match the edited state.
$ws.Activate() | Out-Null
$ws.Application.ActiveWindow.ScrollColumn = 13
$ws.Range("S1").Select() | Out-Null
